$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - copy formatting (bold, border, centered) from G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# "Save" indicator column values (H2:H14)
$values = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 8).Value = $values[$row]
}
